# Weekly price update: insert two new rows (Primera/Segunda) for the
# Brócoli - Terminal La Palmera de La Serena dataset, pushing all
# existing rows down by 2 (dimension grows from A1:R1266 to A1:R1268).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current first data block
# (row 1169), shifting rows 1169:1266 down to 1171:1268.
$ws.Rows("1169:1170").Insert()

# New row 1169 - "Primera" quality, date 45223 (2023-10-24)
$ws.Cells.Item(1169, 1).Value2 = 8
$ws.Cells.Item(1169, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1169, 3).Value2 = "Coquimbo"
$ws.Cells.Item(1169, 4).Value2 = 45223
$ws.Cells.Item(1169, 5).Value2 = 4
$ws.Cells.Item(1169, 6).Value2 = 100112023
$ws.Cells.Item(1169, 7).Value2 = "Brócoli"
$ws.Cells.Item(1169, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1169, 9).Value2 = "Primera"
$ws.Cells.Item(1169, 10).Value2 = 2000
$ws.Cells.Item(1169, 11).Value2 = 800
$ws.Cells.Item(1169, 12).Value2 = 900
$ws.Cells.Item(1169, 13).Value2 = 850
$ws.Cells.Item(1169, 14).Value2 = "`$/unidad"
$ws.Cells.Item(1169, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(1169, 16).Value2 = 850
$ws.Cells.Item(1169, 17).Value2 = 1
$ws.Cells.Item(1169, 18).Value2 = "Hortaliza"

# New row 1170 - "Segunda" quality, date 45223 (2023-10-24)
$ws.Cells.Item(1170, 1).Value2 = 8
$ws.Cells.Item(1170, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1170, 3).Value2 = "Coquimbo"
$ws.Cells.Item(1170, 4).Value2 = 45223
$ws.Cells.Item(1170, 5).Value2 = 4
$ws.Cells.Item(1170, 6).Value2 = 100112023
$ws.Cells.Item(1170, 7).Value2 = "Brócoli"
$ws.Cells.Item(1170, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1170, 9).Value2 = "Segunda"
$ws.Cells.Item(1170, 10).Value2 = 1200
$ws.Cells.Item(1170, 11).Value2 = 600
$ws.Cells.Item(1170, 12).Value2 = 700
$ws.Cells.Item(1170, 13).Value2 = 650
$ws.Cells.Item(1170, 14).Value2 = "`$/unidad"
$ws.Cells.Item(1170, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(1170, 16).Value2 = 650
$ws.Cells.Item(1170, 17).Value2 = 1
$ws.Cells.Item(1170, 18).Value2 = "Hortaliza"
